$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "'" + '27.911.26'
$ws.Range("E2").Value = '  +1.14%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "'" + '1.638.06'

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "'" + '213.62'
$ws.Range("E5").Value = '  +0.66%  '

# Row 6
$ws.Range("E6").Value = '  +0.68%  '

# Row 7
$ws.Range("E7").Value = '  +0.14%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "'" + '23.59'
$ws.Range("E8").Value = '  +0.91%  '

# Row 9
$ws.Range("E9").Value = '  -1.07%  '

# Row 10
$ws.Range("E10").Value = '  +0.61%  '

# Row 11
$ws.Range("E11").Value = '  -0.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "'" + '1.870.52'
$ws.Range("E12").Value = '  +0.55%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "'" + '1.638.20'
$ws.Range("E13").Value = '  +0.53%  '

# Row 14
$ws.Range("E14").Value = '  +4.21%  '

# Row 15
$ws.Range("E15").Value = '  +1.11%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "'" + '65.98'
$ws.Range("E16").Value = '  +0.89%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "'" + '27.906.42'
$ws.Range("E17").Value = '  +1.24%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "'" + '231.76'
$ws.Range("E18").Value = '  +0.28%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "'" + '0.0' + [char]0x2083 + '0723'
$ws.Range("E19").Value = '  +0.63%  '

# Row 20
$ws.Range("E20").Value = '  +0.49%  '

# Row 21
$ws.Range("E21").Value = '  +0.08%  '

# Row 22
$ws.Range("E22").Value = '  +1.90%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "'" + '4.35'
$ws.Range("E23").Value = '  -0.05%  '

# Row 24
$ws.Range("E24").Value = '  -3.55%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "'" + '151.70'

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "'" + '6.91'
$ws.Range("E26").Value = '  +0.09%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "'" + '15.72'
$ws.Range("E27").Value = '  +1.16%  '

# Row 28
$ws.Range("E28").Value = '  -0.05%  '

# Row 29
$ws.Range("E29").Value = '  +0.14%  '

# Row 30
$ws.Range("E30").Value = '  +0.99%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "'" + '0.0483'
$ws.Range("E31").Value = '  -0.14%  '

# Row 32
$ws.Range("E32").Value = '  +1.65%  '

# Row 33
$ws.Range("E33").Value = '  +1.44%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "'" + '1.406.70'
$ws.Range("E34").Value = '  -4.29%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "'" + '1.57'
$ws.Range("E35").Value = '  +1.29%  '

# Row 36
$ws.Range("E36").Value = '  +0.66%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "'" + '0.888'
$ws.Range("E37").Value = '  +1.14%  '

# Row 38
$ws.Range("E38").Value = '  +0.09%  '

# Row 39
$ws.Range("B39").Value = 'TrustWalletToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "'" + '0.919'
$ws.Range("E39").Value = '  -1.94%  '

# Row 40
$ws.Range("B40").Value = 'ImmutableX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "'" + '0.554'
$ws.Range("E40").Value = '  -0.26%  '

# Row 41
$ws.Range("E41").Value = '  -0.56%  '

# Row 42
$ws.Range("E42").Value = '  +0.09%  '

# Row 43
$ws.Range("E43").Value = '  +5.17%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "'" + '66.18'
$ws.Range("E44").Value = '  -2.53%  '

# Row 45
$ws.Range("E45").Value = '  +1.38%  '

# Row 46
$ws.Range("E46").Value = '  -0.09%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "'" + '1.779.34'
$ws.Range("E47").Value = '  +0.71%  '

# Row 48
$ws.Range("E48").Value = '  +0.43%  '

# Row 49
$ws.Range("E49").Value = '  +0.80%  '

# Row 50
$ws.Range("E50").Value = '  +0.42%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "'" + '7.61'
$ws.Range("E51").Value = '  -1.49%  '
